# Capstone Progress Update.
# - Add "Best K Nearest Neighbors accuracy score:" summary on knn_classifier_time (row 36)
# - Extend svm_classifier_time with more timing trials + "Support Vector Machine accuracy score:" summary (row 5)
# - Leave view/selection state matching the final save (svm_classifier_time ends up active/selected)

$wb = $excel.ActiveWorkbook

$wsLinReg  = $wb.Worksheets.Item("multiple_linear_regression_time")
$wsKnnReg  = $wb.Worksheets.Item("knn_regressor_time")
$wsKnnCls  = $wb.Worksheets.Item("knn_classifier_time")
$wsSvm     = $wb.Worksheets.Item("svm_classifier_time")

# --- knn_classifier_time: add "Best K Nearest Neighbors accuracy score:" block on row 36 ---
# (written before the svm sheet's new strings so the shared-string table picks up
#  "accuracy" at index 14, matching the upstream edit order)
$wsKnnCls.Range("I36").Value = "K"
$wsKnnCls.Range("J36").Value = "Nearest"
$wsKnnCls.Range("K36").Value = "Neighbors"
$wsKnnCls.Range("L36").Value = "accuracy"
$wsKnnCls.Range("M36").Value = "score:"
$wsKnnCls.Range("N36").Value = 0.61723163841799999

# --- svm_classifier_time: refreshed trial timings (C1:C10) plus many new trials (rows 11-51) ---
$svmTimes = @{
    1 = 11.7909286022186
    2 = 11.873843908309899
    3 = 11.7439749240875
    4 = 11.5901336669921
    5 = 12.1445662975311
    6 = 11.674047708511299
    7 = 11.915800333023
    8 = 11.998715400695801
    9 = 12.008705139160099
    10 = 11.798919916152901
    11 = 12.0826294422149
    12 = 11.876840353012
    13 = 11.8888278007507
    14 = 11.9907236099243
    15 = 11.8888278007507
    16 = 11.811906337738
    17 = 12.009704113006499
    18 = 12.323382616043
    19 = 12.0067074298858
    20 = 12.0536587238311
    21 = 12.058655261993399
    22 = 11.8758416175842
    23 = 11.8428742885589
    24 = 11.947767257690399
    25 = 11.9747400283813
    26 = 11.659062147140499
    27 = 11.7389824390411
    28 = 11.6400814056396
    29 = 11.9567582607269
    30 = 11.5671570301055
    31 = 11.5671570301055
    32 = 11.7859332561492
    33 = 11.557167768478299
    34 = 11.7389817237854
    35 = 11.7709486484527
    36 = 11.714007139205901
    37 = 11.703019142150801
    38 = 11.9108042716979
    39 = 11.790927648544301
    40 = 11.761957645416199
    41 = 11.7829363346099
    42 = 11.6490733623504
    43 = 11.7759437561035
    44 = 11.77894282341
    45 = 11.7140069007873
    46 = 11.7030179500579
    47 = 11.655067205429001
    48 = 11.6850366592407
    49 = 11.680042028427099
    50 = 11.632089853286701
    51 = 12.3970229625701
}

for ($r = 1; $r -le 51; $r++) {
    $wsSvm.Cells.Item($r, 1).Value = "Time"
    $wsSvm.Cells.Item($r, 2).Value = "elapsed:"
    $wsSvm.Cells.Item($r, 3).Value = $svmTimes[$r]
}

# "Support Vector Machine accuracy score:" block on row 5
$wsSvm.Range("H5").Value = "Support"
$wsSvm.Range("I5").Value = "Vector"
$wsSvm.Range("J5").Value = "Machine"
$wsSvm.Range("K5").Value = "accuracy"
$wsSvm.Range("L5").Value = "score:"
$wsSvm.Range("M5").Value = 0.62358757062100001

# --- View/selection state ---
# Select in tab order, ending on svm_classifier_time so it becomes the saved active tab.
$wsLinReg.Range("N26").Select()
$wsKnnReg.Range("H38").Select()
$wsKnnCls.Range("L53").Select()
$wsSvm.Activate()
$wsSvm.Range("M6").Select()
